# Auto update stock data
# Applies the "2025/11/04 -> 2025/11/05" refresh and the accompanying
# metric updates described by the commit diff.

function Set-TextValue($ws, $addr, $val) {
    # Force the cell to remain a text value (the sheet stores every
    # date/metric as text, not as a real number/date) and avoid leaving
    # behind a stray number-format style on the cell.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (AA) ---
Set-TextValue $ws "A2" "2025/11/05"
Set-TextValue $ws "B2" "4.29"

# --- Row 8 (RIO) ---
Set-TextValue $ws "A8" "2025/11/05"
Set-TextValue $ws "B8" "7.39"

# --- Row 14 (NHY) ---
Set-TextValue $ws "A14" "2025/11/05"
Set-TextValue $ws "B14" "2.66"

# --- Row 20 (RS) ---
Set-TextValue $ws "A20" "2025/11/05"
Set-TextValue $ws "B20" "12.23"

# --- Row 26 (KALU) ---
Set-TextValue $ws "A26" "2025/11/05"
Set-TextValue $ws "B26" "9.68"

# --- Row 32 (RYI) ---
Set-TextValue $ws "A32" "2025/11/05"
Set-TextValue $ws "B32" "24.88"

# --- Row 38 (BVB:ALR) ---
Set-TextValue $ws "A38" "2025/11/05"

# --- Row 44 (ULTR) ---
Set-TextValue $ws "A44" "2025/11/05"
Set-TextValue $ws "B44" "11.22"

# --- Row 50 (BHE) ---
Set-TextValue $ws "A50" "2025/11/05"
Set-TextValue $ws "B50" "11.07"
Set-TextValue $ws "C50" "0.29"
Set-TextValue $ws "D50" "4.30"
Set-TextValue $ws "E50" "2.32"
$ws.Range("G50").ClearContents()
$ws.Range("H50").Value = 6

# --- Row 51 (BHE) ---
$ws.Range("G51").ClearContents()
$ws.Range("H51").Value = 6

# --- Row 52 (BHE) ---
$ws.Range("G52").ClearContents()
$ws.Range("H52").Value = 6

# --- Row 53 (BHE) ---
$ws.Range("G53").ClearContents()
$ws.Range("H53").Value = 6

# --- Row 54 (BHE) ---
$ws.Range("G54").ClearContents()
$ws.Range("H54").Value = 6

# --- Row 55 (BHE) ---
$ws.Range("G55").ClearContents()
$ws.Range("H55").Value = 6

# --- Row 56 (CLS) ---
Set-TextValue $ws "A56" "2025/11/05"

# --- Row 62 (JABIL) ---
Set-TextValue $ws "A62" "2025/11/05"
Set-TextValue $ws "B62" "11.38"

# --- Row 68 (FLEX) ---
Set-TextValue $ws "A68" "2025/11/05"
Set-TextValue $ws "B68" "12.88"

# --- Row 74 (MKS) ---
Set-TextValue $ws "A74" "2025/11/05"
Set-TextValue $ws "B74" "14.91"
